$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns B:E to be stored as text so numeric-looking strings
# (e.g. "1.001", "28.962.33") are preserved exactly instead of being
# auto-converted into numbers/dates by Excel.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.962.33'
$ws.Range("E2").Value = '  -1.68%  '
$ws.Range("D3").Value = '1.829.02'
$ws.Range("E3").Value = '  -2.22%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '239.59'
$ws.Range("E5").Value = '  -1.73%  '
$ws.Range("D6").Value = '0.6856'
$ws.Range("E6").Value = '  -3.26%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.07616'
$ws.Range("E8").Value = '  -3.40%  '
$ws.Range("E9").Value = '  -4.69%  '
$ws.Range("D10").Value = '23.36'
$ws.Range("E10").Value = '  -5.52%  '
$ws.Range("D11").Value = '0.07746'
$ws.Range("E11").Value = '  -3.22%  '
$ws.Range("D12").Value = '1.847.57'
$ws.Range("E12").Value = '  -2.17%  '
$ws.Range("D13").Value = '5.042'
$ws.Range("E13").Value = '  -3.61%  '
$ws.Range("D14").Value = '90.20'
$ws.Range("E14").Value = '  -4.20%  '
$ws.Range("D15").Value = '0.6730'
$ws.Range("E15").Value = '  -4.69%  '
$ws.Range("D16").Value = '6.446'
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("D17").Value = '0.000008268'
$ws.Range("E17").Value = '  -1.18%  '
$ws.Range("D18").Value = '28.979.72'
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("D19").Value = '242.55'
$ws.Range("E19").Value = '  -5.62%  '
$ws.Range("D20").Value = '2.102.75'
$ws.Range("E20").Value = '  -1.32%  '
$ws.Range("D21").Value = '12.66'
$ws.Range("E21").Value = '  -4.41%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").Value = '7.429'
$ws.Range("E23").Value = '  -2.73%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").Value = '0.1471'
$ws.Range("E25").Value = '  -5.77%  '
$ws.Range("D26").Value = '161.27'
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("D27").Value = '8.723'
$ws.Range("E27").Value = '  -3.95%  '
$ws.Range("D28").Value = '18.12'
$ws.Range("E28").Value = '  -4.23%  '
$ws.Range("D29").Value = '1.532'
$ws.Range("E29").Value = '  +2.00%  '
$ws.Range("D30").Value = '4.199'
$ws.Range("E30").Value = '  -3.39%  '
$ws.Range("E31").Value = '  -2.31%  '
$ws.Range("E32").Value = '  -1.56%  '
$ws.Range("D33").Value = '0.05115'
$ws.Range("E33").Value = '  -3.91%  '
$ws.Range("D34").Value = '0.7574'
$ws.Range("E34").Value = '  +0.75%  '
$ws.Range("D35").Value = '1.815'
$ws.Range("E35").Value = '  -4.41%  '
$ws.Range("D36").Value = '1.148'
$ws.Range("E36").Value = '  -2.36%  '
$ws.Range("D37").Value = '2.701'
$ws.Range("E37").Value = '  -0.50%  '
$ws.Range("D38").Value = '0.01836'
$ws.Range("E38").Value = '  -2.54%  '
$ws.Range("D39").Value = '1.217.19'
$ws.Range("E39").Value = '  -4.43%  '
$ws.Range("D40").Value = '2.706'
$ws.Range("E40").Value = '  -1.81%  '
$ws.Range("D41").Value = '0.9114'
$ws.Range("E41").Value = '  +1.04%  '
$ws.Range("D42").Value = '108.69'
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '2.001.47'
$ws.Range("E44").Value = '  -1.28%  '
$ws.Range("D45").Value = '5.405'
$ws.Range("E45").Value = '  -9.62%  '
$ws.Range("D46").Value = '0.5176'
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.00000000122'
$ws.Range("E47").Value = '  -6.28%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.444'
$ws.Range("E48").Value = '  -1.13%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '63.33'
$ws.Range("E49").Value = '  -11.69%  '
$ws.Range("D50").Value = '1.726'
$ws.Range("E50").Value = '  -3.82%  '
$ws.Range("D51").Value = '6.902'
$ws.Range("E51").Value = '  -2.57%  '

Write-Output "applied cryptos update"
